# Add a "Skill Description" column (full skill name) after SkillCode column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before current column B (SFIA Level), shifting
# SFIA Level / Keycode / Description one column to the right.
$ws.Range("B1").EntireColumn.Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "Skill Description"

# Map of SkillCode -> full skill name for the new "Skill Description" column.
$skillNames = @{
    "Autonomy"   = "Autonomy"
    "Influence"  = "Influence"
    "Complexity" = "Complexity"
    "Knowledge"  = "Knowledge"
    "MADE"       = "MADE"
    "SUPP"       = "Supplier management"
}

# Determine the last used row on the sheet.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $code = $ws.Cells.Item($r, 1).Value2
    if ($code -ne $null -and $skillNames.ContainsKey($code)) {
        $ws.Cells.Item($r, 2).Value = $skillNames[$code]
    }
}
